# COREINTL_holdings.xlsx update
#  - Bump the "as of" date in the confidential disclaimer from 2021-05-18 to 2021-05-19
#  - Refresh the Weight / Percent Change figures for the three holdings rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect it so the cells can be edited.
$ws.Unprotect()

# Update the disclaimer text (shared string used by A7) with the new date.
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-19 for illustrative purposes only and are subject to change."

# Row 2 (EFA)
$ws.Range("D2").Value = 0.8477140057096595
$ws.Range("E2").Value = -0.008796179944709936

# Row 3 (EEM)
$ws.Range("D3").Value = 0.1522859942903404
$ws.Range("E3").Value = -0.002984517813840681

# Row 4 (Total)
$ws.Range("D4").Value = 0.9999999999999999
$ws.Range("E4").Value = -0.007911145198631009

# Restore sheet protection to match the original workbook's state.
$ws.Protect("D382")
